$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last conversation row (20) - history trimmed from 19 to 18 rows of data
$ws.Rows.Item(20).Delete()

# Header row: "Conversation" -> "Content"
$ws.Cells.Item(1,1).Value2 = "Role"
$ws.Cells.Item(1,2).Value2 = "Content"
$ws.Cells.Item(1,3).Value2 = "Response_Time"

# Data rows 2-19: roleA/roleB relabeled conversation + response times
$ws.Cells.Item(2,1).Value2 = "roleA"
$ws.Cells.Item(2,2).Value2 = "What is your name?"
$ws.Cells.Item(2,3).Value2 = 1.544272661209106

$ws.Cells.Item(3,1).Value2 = "roleA"
$ws.Cells.Item(3,2).Value2 = "My name is Cuong. How can I assist you today?"
$ws.Cells.Item(3,3).Value2 = 0.6893007755279541

$ws.Cells.Item(4,1).Value2 = "roleB"
$ws.Cells.Item(4,2).Value2 = "Hello Cuong! How can I assist you today?"
$ws.Cells.Item(4,3).Value2 = 1.36078405380249

$ws.Cells.Item(5,1).Value2 = "roleA"
$ws.Cells.Item(5,2).Value2 = "It seems like there is a conversation log between two roles, ""roleA"" and ""roleB"". The conversation goes as follows:`nRoleA: ""What is your name?""`nRoleA: ""My name is Cuong. How can I assist you today?""`nRoleB: ""Hello Cuong! How can I assist you today?""`nIs there anything specific you would like to know or discuss about this conversation log?"
$ws.Cells.Item(5,3).Value2 = 1.106013059616089
$ws.Rows.Item(5).AutoFit()

$ws.Cells.Item(6,1).Value2 = "roleB"
$ws.Cells.Item(6,2).Value2 = "It looks like you have provided a conversation log between two roles, ""roleA"" and ""roleB"". The conversation includes RoleA asking for the name, introducing themselves as Cuong, and RoleB responding with a greeting. `nIs there anything specific you would like to know or discuss about this conversation log?"
$ws.Cells.Item(6,3).Value2 = 0
$ws.Rows.Item(6).AutoFit()

$ws.Cells.Item(7,1).Value2 = "Separator"
$ws.Cells.Item(7,2).Value2 = "-------------------"
$ws.Cells.Item(7,3).Value2 = 0

$ws.Cells.Item(8,1).Value2 = "roleB"
$ws.Cells.Item(8,2).Value2 = "What is your name?"
$ws.Cells.Item(8,3).Value2 = 0.67795729637146

$ws.Cells.Item(9,1).Value2 = "roleA"
$ws.Cells.Item(9,2).Value2 = "Hello! My name is Cuong. How can I assist you today?"
$ws.Cells.Item(9,3).Value2 = 0.6001553535461426

$ws.Cells.Item(10,1).Value2 = "roleB"
$ws.Cells.Item(10,2).Value2 = "Hello Cuong! How can I assist you today?"
$ws.Cells.Item(10,3).Value2 = 0.7574994564056396

$ws.Cells.Item(11,1).Value2 = "roleA"
$ws.Cells.Item(11,2).Value2 = "Hello! My name is Cuong. How can I assist you today?"
$ws.Cells.Item(11,3).Value2 = 0.9494054317474365

$ws.Cells.Item(12,1).Value2 = "roleB"
$ws.Cells.Item(12,2).Value2 = "It seems like there is a conversation between two roles, roleA and roleB. RoleA introduces themselves as Cuong and offers assistance, while roleB responds with a greeting and asks how they can assist. Is there anything specific you would like me to do with this conversation data?"
$ws.Cells.Item(12,3).Value2 = 0

$ws.Cells.Item(13,1).Value2 = "Separator"
$ws.Cells.Item(13,2).Value2 = "-------------------"
$ws.Cells.Item(13,3).Value2 = 0

$ws.Cells.Item(14,1).Value2 = "Separator"
$ws.Cells.Item(14,2).Value2 = "-------------------"
$ws.Cells.Item(14,3).Value2 = 0

$ws.Cells.Item(15,1).Value2 = "Separator"
$ws.Cells.Item(15,2).Value2 = "-------------------"
$ws.Cells.Item(15,3).Value2 = 0

$ws.Cells.Item(16,1).Value2 = "Separator"
$ws.Cells.Item(16,2).Value2 = "-------------------"
$ws.Cells.Item(16,3).Value2 = 0

$ws.Cells.Item(17,1).Value2 = "Separator"
$ws.Cells.Item(17,2).Value2 = "-------------------"
$ws.Cells.Item(17,3).Value2 = 0

$ws.Cells.Item(18,1).Value2 = "Separator"
$ws.Cells.Item(18,2).Value2 = "-------------------"
$ws.Cells.Item(18,3).Value2 = 0

$ws.Cells.Item(19,1).Value2 = "Separator"
$ws.Cells.Item(19,2).Value2 = "-------------------"
$ws.Cells.Item(19,3).Value2 = 0

Write-Output "done"